$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6347.4
$ws.Range("I51").Value = 6112.3335
$ws.Range("K51").Value = 6112.3335
$ws.Range("M51").Value = -5628.3335
$ws.Range("H76").Value = 3011
$ws.Range("I76").Value = 2941.5
$ws.Range("J76").Value = 3150
$ws.Range("K76").Value = 2941.5
$ws.Range("L76").Value = 3150
$ws.Range("M76").Value = -2626.5
$ws.Range("N76").Value = -3780
$ws.Range("H79").Value = 3011
$ws.Range("I79").Value = 2941.5
$ws.Range("J79").Value = 3150
$ws.Range("K79").Value = 2941.5
$ws.Range("L79").Value = 3150
$ws.Range("M79").Value = -1849.5
$ws.Range("N79").Value = -5334
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H96").Value = 1060.125
$ws.Range("I96").Value = 1792
$ws.Range("J96").Value = 328.25
$ws.Range("K96").Value = 5376
$ws.Range("L96").Value = 984.75
$ws.Range("M96").Value = -4003
$ws.Range("N96").Value = -3730.75
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H103").Value = 71429660
$ws.Range("I103").Value = 125000860
$ws.Range("J103").Value = 1400
$ws.Range("K103").Value = 375002580
$ws.Range("L103").Value = 4200
$ws.Range("M103").Value = -375001994
$ws.Range("N103").Value = -5372
$ws.Range("H113").Value = 4977
$ws.Range("I113").Value = 3305
$ws.Range("K113").Value = 3305
$ws.Range("M113").Value = -51
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H127").Value = 7696.5
$ws.Range("I127").Value = 10189.667
$ws.Range("K127").Value = 30569.001
$ws.Range("M127").Value = -25609.001
$ws.Range("H135").Value = 2129.625
$ws.Range("I135").Value = 386.75
$ws.Range("K135").Value = 3480.75
$ws.Range("M135").Value = -945.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1171
$ws.Range("I20").Value = 1006.5
$ws.Range("K20").Value = 1006.5
$ws.Range("M20").Value = -759.5
$ws.Range("H38").Value = 9891
$ws.Range("I38").Value = 32
$ws.Range("J38").Value = 19750
$ws.Range("K38").Value = 32
$ws.Range("L38").Value = 19750
$ws.Range("M38").Value = 384
$ws.Range("N38").Value = -20582
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").ClearContents()
$ws.Range("N115").Value = 0
$ws.Range("H134").Value = 3218.75
$ws.Range("I134").Value = 3625.6667
$ws.Range("J134").Value = 1998
$ws.Range("K134").Value = 10877.0001
$ws.Range("L134").Value = 5994
$ws.Range("M134").Value = -8342.000100000001
$ws.Range("N134").Value = -11064

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 5000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H25").Value = 2073.2856
$ws.Range("H31").Value = 4357.4287
$ws.Range("I31").Value = 3364.2727
$ws.Range("J31").Value = 7999
$ws.Range("K31").Value = 3364.2727
$ws.Range("L31").Value = 7999
$ws.Range("M31").Value = -3069.2727
$ws.Range("N31").Value = -8589
$ws.Range("H34").Value = 4357.4287
$ws.Range("I34").Value = 3364.2727
$ws.Range("J34").Value = 7999
$ws.Range("K34").Value = 3364.2727
$ws.Range("L34").Value = 7999
$ws.Range("M34").Value = -3162.2727
$ws.Range("N34").Value = -8403
$ws.Range("H44").Value = 500
$ws.Range("I44").Value = 500
$ws.Range("K44").Value = 500
$ws.Range("M44").Value = -58
$ws.Range("H55").Value = 1000
$ws.Range("I55").Value = 1000
$ws.Range("K55").Value = 1000
$ws.Range("M55").Value = -685
$ws.Range("H58").Value = 2863.4375
$ws.Range("I58").Value = 2348.0833
$ws.Range("K58").Value = 2348.0833
$ws.Range("M58").Value = -2145.0833
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 2339.4119
$ws.Range("I132").Value = 2264.2307
$ws.Range("J132").Value = 2583.75
$ws.Range("K132").Value = 6792.6921
$ws.Range("L132").Value = 7751.25
$ws.Range("M132").Value = -4262.6921
$ws.Range("N132").Value = -12811.25
$ws.Range("H134").Value = 3492.7
$ws.Range("I134").Value = 3492.7
$ws.Range("K134").Value = 10478.1
$ws.Range("M134").Value = -7943.099999999999
$ws.Range("H136").Value = 2863.4375
$ws.Range("I136").Value = 2348.0833
$ws.Range("K136").Value = 7044.249899999999
$ws.Range("M136").Value = -4494.249899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 298.25
$ws.Range("I18").Value = 298.25
$ws.Range("K18").Value = 894.75
$ws.Range("M18").Value = -725.75
$ws.Range("H68").Value = 5500
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 15000
$ws.Range("M68").Value = -14189
$ws.Range("H71").Value = 5500
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 45000
$ws.Range("M71").Value = -40944
$ws.Range("H98").Value = 295.8
$ws.Range("J98").Value = 545
$ws.Range("L98").Value = 1635
$ws.Range("N98").Value = -4631
$ws.Range("H112").Value = 10654.238
$ws.Range("I112").Value = 1842
$ws.Range("J112").Value = 11581.842
$ws.Range("K112").Value = 5526
$ws.Range("L112").Value = 34745.526
$ws.Range("M112").Value = -4418
$ws.Range("N112").Value = -36961.526
$ws.Range("H122").Value = 427.5
$ws.Range("I122").Value = 350
$ws.Range("J122").Value = 505
$ws.Range("K122").Value = 3150
$ws.Range("L122").Value = 4545
$ws.Range("M122").Value = -700
$ws.Range("N122").Value = -9445

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 219116.5
$ws.Range("I3").Value = 253175
$ws.Range("K3").Value = 253175
$ws.Range("M3").Value = -253059
$ws.Range("H42").Value = 51857.2
$ws.Range("J42").Value = 51857.2
$ws.Range("L42").Value = 51857.2
$ws.Range("N42").Value = -52827.2
$ws.Range("H102").Value = 2063.1428
$ws.Range("I102").Value = 2073.6667
$ws.Range("K102").Value = 2073.6667
$ws.Range("M102").Value = -451.6667000000002
$ws.Range("H115").Value = 51857.2
$ws.Range("J115").Value = 51857.2
$ws.Range("L115").Value = 51857.2
$ws.Range("N115").Value = -54207.2
$ws.Range("H126").Value = 800
$ws.Range("I126").Value = 800
$ws.Range("K126").Value = 2400
$ws.Range("M126").Value = 70

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1274.25
$ws.Range("I46").Value = 1274.25
$ws.Range("K46").Value = 1274.25
$ws.Range("M46").Value = -1086.25
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H82").Value = 4960.1
$ws.Range("I82").Value = 3116.5
$ws.Range("J82").Value = 7725.5
$ws.Range("K82").Value = 3116.5
$ws.Range("L82").Value = 7725.5
$ws.Range("M82").Value = -2755.5
$ws.Range("N82").Value = -8447.5
$ws.Range("H85").Value = 4960.1
$ws.Range("I85").Value = 3116.5
$ws.Range("J85").Value = 7725.5
$ws.Range("K85").Value = 3116.5
$ws.Range("L85").Value = 7725.5
$ws.Range("M85").Value = -1868.5
$ws.Range("N85").Value = -10221.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1967.1052
$ws.Range("I81").Value = 1828.5
$ws.Range("J81").Value = 2706.3333
$ws.Range("K81").Value = 3657
$ws.Range("L81").Value = 5412.6666
$ws.Range("M81").Value = -2596
$ws.Range("N81").Value = -7534.6666
$ws.Range("H84").Value = 1967.1052
$ws.Range("I84").Value = 1828.5
$ws.Range("J84").Value = 2706.3333
$ws.Range("K84").Value = 18285
$ws.Range("L84").Value = 27063.333
$ws.Range("M84").Value = -12981
$ws.Range("N84").Value = -37671.333
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
$ws.Range("H122").Value = 3313.6667
$ws.Range("I122").Value = 3313.6667
$ws.Range("K122").Value = 9941.000100000001
$ws.Range("M122").Value = -7491.000100000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("N130").Value = 0
